# Apply cryptos list update (Mon Mar 11 14:38:44 UTC 2024)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet


# Row 2
$ws.Range("D2").Value = "'72.090.27"
$ws.Range("E2").Value = "  +4.07%  "

# Row 3
$ws.Range("D3").Value = "'4.045.43"
$ws.Range("E3").Value = "  +3.64%  "

# Row 4
$ws.Range("E4").Value = "  +0.09%  "

# Row 5
$ws.Range("D5").Value = "'518.32"

# Row 6
$ws.Range("D6").Value = "'148.75"
$ws.Range("E6").Value = "  +3.61%  "

# Row 7
$ws.Range("D7").Value = "'0.621"
$ws.Range("E7").Value = "  +1.82%  "

# Row 8
$ws.Range("E8").Value = "  +0.19%  "

# Row 9
$ws.Range("D9").Value = "'0.737"
$ws.Range("E9").Value = "  +2.24%  "

# Row 10
$ws.Range("D10").Value = "'0.173"
$ws.Range("E10").Value = "  +2.32%  "

# Row 11
$ws.Range("D11").Value = "'0.0000332"
$ws.Range("E11").Value = "  +0.63%  "

# Row 12
$ws.Range("D12").Value = "'48.15"
$ws.Range("E12").Value = "  +14.69%  "

# Row 13
$ws.Range("E13").Value = "  +6.39%  "

# Row 14
$ws.Range("D14").Value = "'4.687.26"
$ws.Range("E14").Value = "  +3.53%  "

# Row 15
$ws.Range("D15").Value = "'4.033.57"
$ws.Range("E15").Value = "  +2.87%  "

# Row 16
$ws.Range("D16").Value = "'21.18"
$ws.Range("E16").Value = "  +7.68%  "

# Row 17
$ws.Range("D17").Value = "'14.13"
$ws.Range("E17").Value = "  +2.91%  "

# Row 18
$ws.Range("E18").Value = "  -0.37%  "

# Row 19
$ws.Range("E19").Value = "  -2.38%  "

# Row 20
$ws.Range("D20").Value = "'72.136.37"
$ws.Range("E20").Value = "  +4.27%  "

# Row 21
$ws.Range("D21").Value = "'436.39"
$ws.Range("E21").Value = "  +2.15%  "

# Row 22
$ws.Range("D22").Value = "'96.34"
$ws.Range("E22").Value = "  +9.86%  "

# Row 23
$ws.Range("D23").Value = "'3.52"
$ws.Range("E23").Value = "  +5.92%  "

# Row 24
$ws.Range("D24").Value = "'14.63"
$ws.Range("E24").Value = "  +3.36%  "

# Row 25
$ws.Range("B25").Value = "PancakeSwap"
$ws.Range("C25").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D25").Value = "'4.01"
$ws.Range("E25").Value = "  -0.68%  "

# Row 26
$ws.Range("B26").Value = "RenderToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D26").Value = "'11.89"
$ws.Range("E26").Value = "  +3.85%  "

# Row 27
$ws.Range("D27").Value = "'11.18"
$ws.Range("E27").Value = "  +5.76%  "

# Row 28
$ws.Range("D28").Value = "'36.86"
$ws.Range("E28").Value = "  +2.64%  "

# Row 29
$ws.Range("D29").Value = "'3.09"
$ws.Range("E29").Value = "  +10.43%  "

# Row 30
$ws.Range("D30").Value = "'701.20"
$ws.Range("E30").Value = "  +1.09%  "

# Row 31
$ws.Range("D31").Value = "'13.50"
$ws.Range("E31").Value = "  +3.31%  "

# Row 32
$ws.Range("E32").Value = "  +3.21%  "

# Row 33
$ws.Range("D33").Value = "'7.01"
$ws.Range("E33").Value = "  +19.03%  "

# Row 34
$ws.Range("D34").Value = "'68.01"
$ws.Range("E34").Value = "  -1.09%  "

# Row 35
$ws.Range("D35").Value = "'0.0₃0893"
$ws.Range("E35").Value = "  +7.25%  "

# Row 36
$ws.Range("B36").Value = "TheGraph"
$ws.Range("C36").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D36").Value = "'0.435"
$ws.Range("E36").Value = "  -2.65%  "

# Row 37
$ws.Range("B37").Value = "ThetaToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D37").Value = "'3.66"
$ws.Range("E37").Value = "  +24.55%  "

# Row 38
$ws.Range("D38").Value = "'40.42"
$ws.Range("E38").Value = "  +1.11%  "

# Row 39
$ws.Range("D39").Value = "'0.153"
$ws.Range("E39").Value = "  +4.83%  "

# Row 40
$ws.Range("D40").Value = "'0.998"
$ws.Range("E40").Value = "  -0.02%  "

# Row 41
$ws.Range("E41").Value = "  -0.13%  "

# Row 42
$ws.Range("D42").Value = "'0.0487"
$ws.Range("E42").Value = "  +2.01%  "

# Row 43
$ws.Range("E43").Value = "  +2.93%  "

# Row 44
$ws.Range("E44").Value = "  +0.20%  "

# Row 45
$ws.Range("D45").Value = "'3.52"
$ws.Range("E45").Value = "  +5.21%  "

# Row 46
$ws.Range("E46").Value = "  +3.46%  "

# Row 47
$ws.Range("D47").Value = "'3.13"
$ws.Range("E47").Value = "  +3.38%  "

# Row 48
$ws.Range("D48").Value = "'9.04"
$ws.Range("E48").Value = "  +9.04%  "

# Row 49
$ws.Range("B49").Value = "LidoDAOToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D49").Value = "'3.33"
$ws.Range("E49").Value = "  +1.89%  "

# Row 50
$ws.Range("B50").Value = "FLOKI"
$ws.Range("C50").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D50").Value = "'0.000271"
$ws.Range("E50").Value = "  +21.25%  "

# Row 51
$ws.Range("D51").Value = "'0.0₆0341"
$ws.Range("E51").Value = "  +4.38%  "
